# resume.docx edit: bump the "Experience" heading to 12pt with tight
# character spacing, and re-title the ATLAS job entry.

$d = $word.ActiveDocument

# Track changes would wrap our edits in <w:ins>/<w:del>; we want a clean
# in-place edit, matching the target OOXML.
$d.TrackRevisions = $false

# --- 1. "Experience" section heading: spacing -1 (i.e. -0.05pt in COM's
#        point-based Font.Spacing), 12pt (24 half-points) size for both
#        the ascii and complex-script (Cs) size fields. Setting this on
#        the paragraph's own Range (which spans the run AND the trailing
#        paragraph mark) pushes the size values onto both the run rPr and
#        the paragraph-mark rPr in pPr, matching the diff.
$headingFound = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Experience") {
        $r = $para.Range
        $r.Font.Spacing = -0.05
        $r.Font.Size = 12
        $r.Font.SizeBi = 12
        $headingFound = $true
        break
    }
}

# --- 2. Re-title the job entry from "ATLAS Collaboration" to
#        "SMU ATLAS Experiment".
$d.Content.Find.Execute("ATLAS Collaboration", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SMU ATLAS Experiment", 2)

Write-Host "headingFound=$headingFound"
